$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("N-wni")

# Header for new column E
$ws.Range("E1").Value = "Odsetek"

# Row 2 gets its own (non-shared) formula
$ws.Range("E2").Formula = "=C2/D2*100"

# Rows 3 through 41 share one formula definition
$ws.Range("E3:E41").Formula = "=C3/D3*100"

# Update the view: scroll position and selection to match the edited workbook
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("E2:E41").Select()
